$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "harvester" value from "BROWN" to "H.BROWN" for all data rows (2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
}

# Update the active selection to match the saved worksheet view.
[void]$ws.Range("B3:B27").Select()
